$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "to create meaningful phrases" -> "to create phrases"
#   (drop the word "meaningful")
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "to create meaningful phrases", $true, $false, $false, $false, $false,
    $true, 1, $false, "to create phrases", 2)

# ---------------------------------------------------------------------------
# Change 2: "...useful for the user to quickly navigate..."
#        -> "...useful for the learner to quickly navigate..."
#   In the real document this was produced by selecting the word "user" and
#   retyping "learner", which splits the run in three (text before / the new
#   word / text after), with the newly-typed word carrying its own run
#   properties.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "useful for the user to quickly", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

$searchStart = $rng.Start
$wordStart = $searchStart + 15   # offset of "user" within the matched phrase
$wordEnd = $searchStart + 19     # "user" is 4 characters long

$targetRange = $d.Range($wordStart, $wordEnd)
$targetRange.Text = "learner"

$newWordRange = $d.Range($wordStart, $wordStart + 7)   # "learner" is 7 characters long
$newWordRange.Font.Size = 12
$newWordRange.Font.Color = -16777216

# ---------------------------------------------------------------------------
# Change 3: add ", depending on current capability of EducationalWeb" before
# the closing parenthesis.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "this is an optional goal for this project)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "this is an optional goal for this project, depending on current capability of EducationalWeb)",
    2)
